$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data per latest scrape.
# Cells whose new text would otherwise be auto-parsed by Excel as a
# number (e.g. "1.002", "0.3101") are forced to Text format first so
# the literal string is preserved exactly, matching the source data.

$ws.Range('D2').Value = '29.827.01'
$ws.Range('E2').Value = '  -0.52%  '
$ws.Range('D3').Value = '1.865.32'
$ws.Range('E3').Value = '  -1.64%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7233'
$ws.Range('E5').Value = '  -6.32%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '241.69'
$ws.Range('E6').Value = '  -1.20%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.002'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3101'
$ws.Range('E8').Value = '  -1.44%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07035'
$ws.Range('E9').Value = '  -3.55%  '
$ws.Range('B10').Value = 'Solana'
$ws.Range('C10').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.36'
$ws.Range('E10').Value = '  -5.66%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08405'
$ws.Range('E11').Value = '  +4.41%  '
$ws.Range('D12').Value = '1.913.40'
$ws.Range('E12').Value = '  +2.91%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7423'
$ws.Range('E13').Value = '  -4.20%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.306'
$ws.Range('E14').Value = '  -3.52%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.02'
$ws.Range('E15').Value = '  -3.47%  '
$ws.Range('D16').Value = '29.856.39'
$ws.Range('E16').Value = '  +0.16%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.011'
$ws.Range('E17').Value = '  -3.30%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.47'
$ws.Range('E18').Value = '  -4.07%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '239.61'
$ws.Range('E19').Value = '  -3.31%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007771'
$ws.Range('E20').Value = '  -0.95%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.002'
$ws.Range('E21').Value = '  +0.17%  '
$ws.Range('D22').Value = '2.134.49'
$ws.Range('E22').Value = '  +5.21%  '
$ws.Range('E23').Value = '  +0.21%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.869'
$ws.Range('E24').Value = '  -3.47%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1530'
$ws.Range('E25').Value = '  -2.68%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.264'
$ws.Range('E26').Value = '  -2.45%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '162.64'
$ws.Range('E27').Value = '  -0.16%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.48'
$ws.Range('E28').Value = '  -1.71%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.998'
$ws.Range('E29').Value = '  -1.82%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.485'
$ws.Range('E30').Value = '  +4.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.529'
$ws.Range('E31').Value = '  -1.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.483'
$ws.Range('E32').Value = '  -0.90%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.217'
$ws.Range('E33').Value = '  +2.84%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05298'
$ws.Range('E34').Value = '  -4.29%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.225'
$ws.Range('E35').Value = '  -1.94%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7410'
$ws.Range('E36').Value = '  -1.44%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.001'
$ws.Range('E37').Value = '  +0.45%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.712'
$ws.Range('E38').Value = '  +1.06%  '
$ws.Range('E39').Value = '  -0.23%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.747'
$ws.Range('E40').Value = '  -1.67%  '
$ws.Range('D41').Value = '1.114.36'
$ws.Range('E41').Value = '  +2.44%  '
$ws.Range('E42').Value = '  -1.78%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.012'
$ws.Range('E43').Value = '  -0.37%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '71.95'
$ws.Range('E44').Value = '  -3.29%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8669'
$ws.Range('E45').Value = '  +1.60%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.002'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '102.24'
$ws.Range('E47').Value = '  -0.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.607'
$ws.Range('E48').Value = '  +0.03%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.821'
$ws.Range('E49').Value = '  -3.96%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.023'
$ws.Range('E50').Value = '  +0.74%  '
$ws.Range('D51').Value = '2.024.23'
$ws.Range('E51').Value = '  +2.15%  '
